# Redeem points for phone 71277628 (76.0 points) — adds a new
# redemption log row and normalizes the previous row's phone value
# from a text-typed number to a genuine numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11's phone column was stored as text ("71277628"); the app later
# rewrote it as a true number once the redemption settled.
$ws.Range("A11").Value = 71277628

# New redemption row: phone stays text-typed (matches how freshly
# logged rows are written), points = 76, timestamp of the event.
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "71277628"
$ws.Range("A12").ClearFormats()

$ws.Range("B12").Value = 76
$ws.Range("C12").Value = "2025-08-18T16:54:45"
